$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Topic column (D) for rows that have one, in row order except D24 (RMarkdown)
# which is entered last -- this reproduces the shared-string insertion order
# recorded by the original author.
$ws.Range("D23").Value = "Software setup"
$ws.Range("D25").Value = "Advanced programming tips"
$ws.Range("D26").Value = "The shell"
$ws.Range("D27").Value = "Pandas"
$ws.Range("D29").Value = "RMarkdown websites"

# Lab labels (A column)
$ws.Range("A23").Value = "lab01"
$ws.Range("A24").Value = "lab02"
$ws.Range("A25").Value = "lab03"
$ws.Range("A26").Value = "lab04"
$ws.Range("A27").Value = "lab05"
$ws.Range("A28").Value = "lab06"
$ws.Range("A29").Value = "lab07"
$ws.Range("A30").Value = "lab08"
$ws.Range("A31").Value = "lab09"
$ws.Range("A32").Value = "lab10"

# Last new shared string
$ws.Range("D24").Value = "RMarkdown"

# Dates (B column) -- first lab session is a literal date, the rest are
# formulas that add 7 days to the previous week's date.
$ws.Range("B23").Value = 42641
$ws.Range("B24").Formula = "=B23+7"
$ws.Range("B25").Formula = "=B24+7"
$ws.Range("B26").Formula = "=B25+7"
$ws.Range("B27").Formula = "=B26+7"
$ws.Range("B28").Formula = "=B27+7"
$ws.Range("B29").Formula = "=B28+7"
$ws.Range("B30").Formula = "=B29+7"
$ws.Range("B31").Formula = "=B30+7"
$ws.Range("B32").Formula = "=B31+7"

# Completed flags (C column) -- all new lab sessions are not yet completed.
$ws.Range("C23").Value = $false
$ws.Range("C24").Value = $false
$ws.Range("C25").Value = $false
$ws.Range("C26").Value = $false
$ws.Range("C27").Value = $false
$ws.Range("C28").Value = $false
$ws.Range("C29").Value = $false
$ws.Range("C30").Value = $false
$ws.Range("C31").Value = $false
$ws.Range("C32").Value = $false

# Update selection to match the committed state (column E selected)
$ws.Columns("E").Select()
